$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" positioned right before "总计"
# ------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q1"

# NOTE: after Add(before), the variable used as the "before" anchor no
# longer points at the original "总计" sheet (it now aliases the newly
# inserted sheet's slot) - so re-fetch "总计" by name to get a handle
# on the real data sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row for the new "2022-Q1" fund-holdings sheet
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows for "2022-Q1" - fund code / 规模 / 仓位 / 占比 / 市值 are kept
# as text (matching the other quarter sheets), only the rank column is
# a real number.
$data = @(
    @("164811", "工银瑞信中证京津冀协同发展主题指数（LOF）A", "0.23", "94.28", "4.18", "0.0096", 2),
    @("512780", "广发中证京津冀协同发展主题ETF", "0.13", "98.52", "3.29", "0.0043", 4),
    @("164825", "工银瑞信中证京津冀协同发展主题指数（LOF）C", "0.06", "94.28", "4.18", "0.0025", 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $newSheet.Range("A$row").Value = $i

    $newSheet.Range("B$row").NumberFormat = "@"
    $newSheet.Range("B$row").Value = $item[0]
    $newSheet.Range("B$row").ClearFormats()

    $newSheet.Range("C$row").Value = $item[1]

    $newSheet.Range("D$row").NumberFormat = "@"
    $newSheet.Range("D$row").Value = $item[2]
    $newSheet.Range("D$row").ClearFormats()

    $newSheet.Range("E$row").NumberFormat = "@"
    $newSheet.Range("E$row").Value = $item[3]
    $newSheet.Range("E$row").ClearFormats()

    $newSheet.Range("F$row").NumberFormat = "@"
    $newSheet.Range("F$row").Value = $item[4]
    $newSheet.Range("F$row").ClearFormats()

    $newSheet.Range("G$row").NumberFormat = "@"
    $newSheet.Range("G$row").Value = $item[5]
    $newSheet.Range("G$row").ClearFormats()

    $newSheet.Range("H$row").Value = $item[6]

    $idxCell = $newSheet.Range("A$row")
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
}

# ------------------------------------------------------------------
# 2) Update "总计" sheet: insert a new top data row summarizing
#    the 2022-Q1 quarter, pushing existing rows down.
# ------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.02

# Renumber the index column for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Match formatting used by the sibling rows
$totalIdxCell = $totalSheet.Range("A2")
$totalIdxCell.Font.Bold = $true
$totalIdxCell.Borders.LineStyle = 1
$totalIdxCell.HorizontalAlignment = -4108
$totalIdxCell.VerticalAlignment = -4160

$totalSheet.Range("B2:D2").ClearFormats()

# Restore the originally active sheet/tab (Worksheets.Add() activates the
# newly created sheet as a side effect).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "2022-Q1 data added"
